# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Row 3: replace the top "Bad Driver" entry with new driver/stats
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.200.0.6"
$ws.Range("C3").Value = 155
$ws.Range("D3").Value = 98.90000000000001

# Row 4: Totals row - Critical Minutes total updates to match new C3 value
$ws.Range("C4").Value = 155

# Row 12 & 13: Total Samples counts updated for good drivers
$ws.Range("B12").Value = 11140
$ws.Range("B13").Value = 14487
